# Generate Report for Handoff
# Updates the row for "cae07396-fd5e-4559-ae4b-37b13a9513e0" (the e2e file that is
# now ready for handoff) across the Overview, zh-cn and de-de sheets, and widens
# the "Error Detail" column to fit the new error message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad9862d0bde3549ec95c88f2e7455dd4f9f030a1/e2e/cae07396-fd5e-4559-ae4b-37b13a9513e0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c3880711b4cbb0e68c1ecbaee7cd66fcb1ae5744/e2e/cae07396-fd5e-4559-ae4b-37b13a9513e0.md."

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-17 14:48:52"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-17 14:48:47"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-17 14:48:52"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15

Write-Host "Generate Report for Handoff: applied"
